$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ("Ping"): status flips from "Incomplete" to "In Progress"
#     now that the ping command accounts for CSD/Bot Permissions, taking
#     on the yellow "In Progress" styling that F3 currently has (copy it
#     first, before F3's own style changes below).
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F7").Value = "In Progress"

# --- Row 3 ("Joke"): status flips from "In Progress" to a brand new
#     "Complete 0.3.1.2b" status, picking up the green "Complete" styling
#     already used by other complete rows (e.g. F24), and a new Line
#     Complete number shows up in G3.
$ws.Range("F24").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F3").Value = "Complete 0.3.1.2b"
$ws.Range("G3").Value = 156

# --- Move the active selection to match the author's final cursor spot.
$ws.Range("F9").Select()
